$d = $word.ActiveDocument

$replacements = @(
    @("53÷7=7, 4", "17÷6=2, 5"),
    @("56÷4=14, 0", "92÷5=18, 2"),
    @("19÷7=2, 5", "66÷3=22, 0"),
    @("86÷7=12, 2", "95÷4=23, 3"),
    @("91÷7=13, 0", "55÷5=11, 0"),
    @("92÷7=13, 1", "91÷6=15, 1"),
    @("85÷2=42, 1", "54÷3=18, 0"),
    @("53÷4=13, 1", "49÷4=12, 1"),
    @("96÷8=12, 0", "79÷7=11, 2"),
    @("23÷3=7, 2", "84÷5=16, 4"),
    @("20÷9=2, 2", "29÷6=4, 5"),
    @("96÷2=48, 0", "37÷8=4, 5"),
    @("58÷3=19, 1", "26÷5=5, 1"),
    @("94÷2=47, 0", "11÷8=1, 3"),
    @("19÷2=9, 1", "48÷3=16, 0"),
    @("85÷6=14, 1", "55÷8=6, 7"),
    @("29÷9=3, 2", "48÷7=6, 6"),
    @("15÷6=2, 3", "36÷5=7, 1"),
    @("67÷8=8, 3", "62÷4=15, 2"),
    @("30÷5=6, 0", "15÷9=1, 6"),
    @("23÷2=11, 1", "71÷7=10, 1"),
    @("77÷3=25, 2", "90÷4=22, 2"),
    @("51÷8=6, 3", "97÷5=19, 2"),
    @("94÷9=10, 4", "80÷7=11, 3"),
    @("41÷6=6, 5", "20÷5=4, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
